$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Update row 2: A2 becomes "my a", B2 stays 1, C2 becomes 7
$ws.Range("A2").Value = "my a"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 7

# Remove rows 3 through 5 (they are no longer part of the data)
$ws.Range("A3:C5").EntireRow.Delete()
